$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    'ABCB1',
    'ADCY9',
    'AK4',
    'ARG2',
    'ATXN1',
    'BMP6',
    'BUB1',
    'C1QA',
    'C2',
    'CASP4',
    'COL3A1',
    'DLAT',
    'DLGAP5',
    'EIF2AK1',
    'EPHB3',
    'ETFDH',
    'FGF13',
    'GPR19',
    'GPRC5C',
    'HBEGF',
    'HIF1A',
    'HS3ST1',
    'IDO1',
    'IL2RA',
    'ISG20',
    'KRT1',
    'LRP1',
    'MAP2K5',
    'MICB',
    'MMP9',
    'MYBPC3',
    'NQO2',
    'PARP12',
    'PDHB',
    'PPA1',
    'PTGDS',
    'PTGER3',
    'PTPRJ',
    'RABGAP1L',
    'RNF123',
    'RSAD2',
    'SCD',
    'SLC6A6',
    'SLC7A1',
    'UGDH',
    'WEE1',
    'CR2',
    'MICAL3',
    'PRSS21',
    'RPH3A',
    'AKAP9',
    'CACNA1I',
    'CD1E',
    'CD244',
    'DHRS9',
    'DYNC1LI1',
    'EIF5A',
    'ERV3-2',
    'FCGBP',
    'GRAMD1C',
    'GSTM1',
    'HERC2',
    'HLA-DPB1',
    'HLA-DRB4',
    'HOXB2',
    'JAK3',
    'KLRG1',
    'LAPTM4B',
    'LOC389906',
    'MOCS1',
    'NEFH',
    'PAM',
    'PRKD3',
    'RPH3AL',
    'SCAND2P',
    'SERPINE2',
    'STARD8',
    'SUPT3H',
    'SUSD4',
    'TRIM2',
    'ALAD',
    'AMFR',
    'ANO1',
    'APAF1',
    'ATP1B1',
    'BIN1',
    'C3',
    'CBR3',
    'CDK1',
    'CDKN2C',
    'CES1',
    'CHEK1',
    'CHPT1',
    'CSAD',
    'CXCL6',
    'CYP27A1',
    'DAPK2',
    'DUSP5',
    'E2F8',
    'EFHC1',
    'EPHA5',
    'FZD2',
    'GINS3',
    'HMMR',
    'HTATIP2',
    'IFI44L',
    'IFNAR2',
    'ISCA1',
    'ITGB4',
    'KIF11',
    'LCP2',
    'LHX2',
    'MAOA',
    'MAST3',
    'MGST3',
    'MTMR10',
    'NDUFS1',
    'NINJ1',
    'NR4A2',
    'OAS1',
    'PAX4',
    'PECAM1',
    'PLPP3',
    'PPFIA4',
    'PTPRD',
    'PTPRE',
    'PTTG1',
    'PYGM',
    'RAP2B',
    'RYR2',
    'SLC1A5',
    'SLC22A18',
    'SLC23A2',
    'SPP1',
    'STAT1',
    'SYNCRIP',
    'TIMM13',
    'TLR6',
    'TMBIM6',
    'TNFAIP6',
    'TOP2A',
    'TPST1',
    'TSTA3',
    'ZFPM2',
    'PPP4R1',
    'RAB20',
    'AVIL',
    'BBS9',
    'CEP63',
    'CYP2C9',
    'FTH1',
    'HNRNPL',
    'KRT23',
    'PCYT1A',
    'RERE',
    'RNF24',
    'TNFRSF10B',
    'ABLIM1',
    'ACADM',
    'AMD1',
    'BHLHE40',
    'CBX5',
    'CCND2',
    'CD3G',
    'CFLAR',
    'CKB',
    'CR1',
    'DDX18',
    'DEGS1',
    'EPHX2',
    'FH',
    'GAPVD1',
    'GLA',
    'GPX3',
    'HADHA',
    'HBB',
    'HDAC9',
    'HNRNPD',
    'HNRNPU',
    'HSP90AB1',
    'HSPA4',
    'HSPA9',
    'HSPH1',
    'IGSF6',
    'KPNB1',
    'LBR',
    'NUMA1',
    'NUSAP1',
    'OASL',
    'PPP1R15A',
    'RAB31',
    'RAD21',
    'RBPJ',
    'SLC16A1',
    'SNX10',
    'SPOCK2',
    'STAT2',
    'SUCLA2',
    'TES',
    'TFRC',
    'TIMM17A',
    'WAS',
    'WDR43',
    'CC2D1A',
    'DHX34',
    'EMC3',
    'EXT1',
    'HLA-DQB1',
    'LPAR1',
    'MGAM',
    'MPZL1',
    'PADI4',
    'SERPINF1',
    'API5',
    'CEACAM1',
    'CHMP1B',
    'DPYSL2',
    'EIF4G1',
    'GIT2',
    'LMO4',
    'OSGIN2'
)

$n = $values.Length
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

Write-Output "Set $n values in column A"
